# Daily attendance processing - 2025-11-12 22:21:50
#
# The "Recorded By" column (G) lists the editors who touched each
# attendance row as a comma-separated string, e.g.
#   "System, dnasr281@gmail.com"
# For every row where "System" (any case) is the first or last name in
# that list, swap the first and last entries so "System" moves to the
# end (and whichever name was last moves to the front). Rows whose list
# doesn't mention "System" at either end (e.g. a plain e-mail list) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $first = $parts[0]
    $last = $parts[$parts.Count - 1]

    if ($first.ToLower() -eq "system" -or $last.ToLower() -eq "system") {
        $tmp = $parts[0]
        $parts[0] = $parts[$parts.Count - 1]
        $parts[$parts.Count - 1] = $tmp
        $cell.Value = [string]::Join(", ", $parts)
    }
}
